# Append the 12/28/2025 profit row (row 34) to Sheet1, matching the
# "Update profit files after running on 2025-12-28" commit.
#
# Column A holds the date as literal text (e.g. "11/26/2025") in every
# existing row, not a real Excel date serial. A leading apostrophe forces
# the new cell to be entered as text too, instead of being auto-converted
# into a date value. The apostrophe itself becomes the quote-prefix marker
# and is not part of the stored text.
#
# Style note: typing a quoted value into a General cell normally tags the
# cell with a dedicated "quote prefix" style. Resetting the style back to
# Normal afterwards keeps the cell's text content/type intact while
# dropping that extra style, so A34 ends up using the same default
# (unstyled) formatting as the rest of the data rows, like the other new
# cells B34:J34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Value = "'12/28/2025"
$ws.Range("A34").Style = "Normal"

$ws.Range("B34").Value = 12175.48
$ws.Range("C34").Value = 0.21143141114929
$ws.Range("D34").Value = 0.78856858885071
$ws.Range("E34").Value = -139.07
$ws.Range("F34").Value = -25.99
$ws.Range("G34").Value = -20990.56
$ws.Range("H34").Value = -68.62
$ws.Range("I34").Value = -478.57
$ws.Range("J34").Value = -15.68
